$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.421875
    "B3" = 0.453125
    "B4" = 0.390625
    "B5" = 0.3125
    "B6" = 0.296875
    "B7" = 0.296875
    "B8" = 0.28125
    "B9" = 0.3125
    "B10" = 0.28125
    "B11" = 0.28125
    "B12" = 0.25
    "B14" = 0.34375
    "B15" = 0.453125
    "B16" = 0.375
    "B17" = 0.390625
    "B18" = 0.3125
    "B19" = 0.25
    "B20" = 0.21875
    "B21" = 0.265625
    "B22" = 0.265625
    "B23" = 0.21875
    "B24" = 0.21875
    "B25" = 0.21875
    "B26" = 0.265625
    "B27" = 0.25
    "B28" = 0.234375
    "B29" = 0.28125
    "B30" = 0.25
    "B31" = 0.265625
    "B32" = 0.25
    "B34" = 0.21875
    "B35" = 0.21875
    "B36" = 0.21875
    "B37" = 0.21875
    "B38" = 0.21875
    "B39" = 0.21875
    "B40" = 0.21875
    "B42" = 0.21875
    "B43" = 0.21875
    "B44" = 0.21875
    "B45" = 0.203125
    "B46" = 0.234375
    "B47" = 0.234375
    "B48" = 0.234375
    "B49" = 0.234375
    "B50" = 0.234375
    "B51" = 0.21875
    "B52" = 0.21875
    "B53" = 0.21875
    "B54" = 0.21875
    "B55" = 0.21875
    "B56" = 0.21875
    "B57" = 0.21875
    "B58" = 0.21875
    "B59" = 0.21875
    "B60" = 0.21875
    "B61" = 0.21875
    "B62" = 0.21875
    "B63" = 0.21875
    "B64" = 0.21875
    "B65" = 0.21875
    "B66" = 0.21875
    "B67" = 0.21875
    "B68" = 0.21875
    "B69" = 0.21875
    "B70" = 0.21875
    "B71" = 0.21875
    "B72" = 0.21875
    "B73" = 0.21875
    "B74" = 0.21875
    "B75" = 0.21875
    "B76" = 0.21875
    "B77" = 0.21875
    "B78" = 0.21875
    "B79" = 0.21875
    "B80" = 0.21875
    "B81" = 0.21875
    "B82" = 0.21875
    "B83" = 0.21875
    "B84" = 0.21875
    "B85" = 0.21875
    "B86" = 0.21875
    "B87" = 0.21875
    "B88" = 0.21875
    "B89" = 0.21875
    "B90" = 0.21875
    "B91" = 0.21875
    "B92" = 0.21875
    "B93" = 0.21875
    "B94" = 0.21875
    "B95" = 0.21875
    "B96" = 0.21875
    "B97" = 0.21875
    "B98" = 0.21875
    "B99" = 0.21875
    "B100" = 0.21875
    "B101" = 0.21875
    "B102" = 0.21875
    "B104" = 0.1875
    "B105" = 0.25
    "B106" = 0.1875
    "B107" = 0.3125
    "B108" = 0.234375
    "B109" = 0.171875
    "B110" = 0.21875
    "B111" = 0.234375
    "B112" = 0.15625
    "B113" = 0.265625
    "B114" = 0.1875
    "B115" = 0.234375
    "B116" = 0.15625
    "B117" = 0.328125
    "B118" = 0.2622950819672131
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}

# Update the inline string repr (DisplayOutputs object memory address) in column A for rows 102-118
$newRepr = "<__main__.DisplayOutputs object at 0x7fe5301b0f10>"
for ($r = 102; $r -le 118; $r++) {
    $ws.Range("A" + $r).Value = $newRepr
}
